# Weekly fruit/vegetable price update: a new week's data block (3 rows,
# Fecha = 44642) is inserted at the top of the "Pepino dulce" quality
# breakdown (rows 136-138), pushing all the later weeks down by 3 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 136 so everything that used to
# live at 136.. now lives at 139.. (dimension grows from R173 to R176).
$ws.Rows("136:138").Insert()

# Constant values shared by every data row in this block.
$mercadoId   = 12
$mercado     = "Mapocho Venta Directa de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$categoriaId = 100112043
$categoria   = "Pepino dulce"
$variedad    = "Cultivar IV Región"
$unidad      = "$/bandeja 18 kilos"
$origen      = "Provincia de Limarí"
$kgUnidades  = 18
$clasif      = "Hortaliza"

# New week's rows: r, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$newRows = @(
    @(136, "Especial", 350, 14000, 14000, 14000, 778),
    @(137, "Primera",  330, 12000, 12000, 12000, 667),
    @(138, "Segunda",  280, 10000, 10000, 10000, 556)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = 44642
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $categoriaId
    $ws.Cells.Item($r, 7).Value  = $categoria
    $ws.Cells.Item($r, 8).Value  = $variedad
    $ws.Cells.Item($r, 9).Value  = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
    $ws.Cells.Item($r, 11).Value = $row[3]
    $ws.Cells.Item($r, 12).Value = $row[4]
    $ws.Cells.Item($r, 13).Value = $row[5]
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $row[6]
    $ws.Cells.Item($r, 17).Value = $kgUnidades
    $ws.Cells.Item($r, 18).Value = $clasif
}
